$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($ref in $cells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '321.29'
$ws.Range('E2').Value = '-3.22%'
$ws.Range('D3').Value = '42.83'
$ws.Range('E3').Value = '-6.15%'
$ws.Range('D4').Value = '5.224'
$ws.Range('E4').Value = '-4.49%'
$ws.Range('D5').Value = '0.08230'
$ws.Range('E5').Value = '-3.46%'
$ws.Range('D6').Value = '4.323'
$ws.Range('E6').Value = '-2.70%'
$ws.Range('D7').Value = '1.779'
$ws.Range('E7').Value = '-14.37%'
$ws.Range('D8').Value = '0.9501'
$ws.Range('E8').Value = '-3.86%'
$ws.Range('D9').Value = '0.1123'
$ws.Range('E9').Value = '-3.41%'
$ws.Range('D10').Value = '0.1881'
$ws.Range('E10').Value = '-2.12%'
$ws.Range('D11').Value = '0.09407'
$ws.Range('E11').Value = '-3.48%'
$ws.Range('D12').Value = '0.04631'
$ws.Range('E12').Value = '-1.96%'
$ws.Range('D13').Value = '7.467'
$ws.Range('E13').Value = '-21.36%'
$ws.Range('D14').Value = '0.1059'
$ws.Range('E14').Value = '-0.04%'
$ws.Range('D15').Value = '0.001291'
$ws.Range('E15').Value = '-0.18%'
$ws.Range('D16').Value = '0.005654'
$ws.Range('E16').Value = '-5.86%'
$ws.Range('D17').Value = '3.353'
$ws.Range('E17').Value = '-0.88%'
$ws.Range('D18').Value = '2.522'
$ws.Range('E18').Value = '-0.51%'
$ws.Range('E19').Value = '0.40%'
$ws.Range('D20').Value = '0.1388'
$ws.Range('E20').Value = '1.11%'
$ws.Range('D21').Value = '0.2548'
$ws.Range('E21').Value = '-0.17%'
$ws.Range('D22').Value = '0.04167'
$ws.Range('E22').Value = '0.70%'
$ws.Range('D23').Value = '0.001250'
$ws.Range('E23').Value = '-4.03%'
$ws.Range('D24').Value = '0.004285'
$ws.Range('E24').Value = '-4.06%'
$ws.Range('D25').Value = '0.0001221'
$ws.Range('E25').Value = '-6.29%'
$ws.Range('D26').Value = '0.0002978'
$ws.Range('E26').Value = '-0.31%'
$ws.Range('D38').Value = '0.02667'
$ws.Range('E38').Value = '-3.49%'
$ws.Range('D39').Value = '0.05639'
$ws.Range('E39').Value = '-1.25%'
$ws.Range('D40').Value = '0.008146'
$ws.Range('E40').Value = '3.24%'
$ws.Range('D41').Value = '0.1406'
$ws.Range('E41').Value = '-1.97%'
$ws.Range('D42').Value = '0.006482'
$ws.Range('E42').Value = '-10.80%'
$ws.Range('D43').Value = '0.002150'
$ws.Range('E43').Value = '3.49%'
$ws.Range('D44').Value = '0.007705'
$ws.Range('E44').Value = '-12.69%'
$ws.Range('D45').Value = '0.3486'
$ws.Range('E45').Value = '-1.79%'
$ws.Range('D46').Value = '0.00006756'
$ws.Range('E46').Value = '-3.83%'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').Value = '-0.14%'
$ws.Range('D48').Value = '0.003076'
$ws.Range('E48').Value = '-11.63%'
$ws.Range('D49').Value = '0.004099'
$ws.Range('E49').Value = '15.83%'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').Value = '-0.14%'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').Value = '-0.14%'
